$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added newest games and changed train/test split to 80/20"
# Row 3 (RF, max_depth=None, min_samples_leaf=1, n_estimators=150, max_features=1/3)
# is re-scored against the refreshed dataset under the new 80/20 split.
$ws.Range("G3").Value = 0.1215               # train score
$ws.Range("H3").Value = 0.31630000000000003  # test score
# I3 = ABS(H3-G3) recalculates automatically.

# New holdout score column entry for this row.
$ws.Range("J3").Value = 0.30430000000000001

# Leave the cursor on the newly added cell.
$null = $ws.Range("J3").Select()
